$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.289.34'
$ws.Range("E2").Value = '  +0.98%  '

$ws.Range("D3").Value = '2.271.13'
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '''305.67'
$ws.Range("E5").Value = '  +0.59%  '

$ws.Range("D6").Value = '''97.73'
$ws.Range("E6").Value = '  +4.84%  '

$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").Value = '''0.494'
$ws.Range("E9").Value = '  +1.64%  '

$ws.Range("D10").Value = '''35.78'
$ws.Range("E10").Value = '  +9.29%  '

$ws.Range("D11").Value = '''0.0797'
$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("E12").Value = '  -1.03%  '

$ws.Range("D13").Value = '''6.67'
$ws.Range("E13").Value = '  -0.57%  '

$ws.Range("D14").Value = '2.597.16'
$ws.Range("E14").Value = '  -1.02%  '

$ws.Range("E15").Value = '  +0.71%  '

$ws.Range("D16").Value = '2.271.61'
$ws.Range("E16").Value = '  +0.02%  '

$ws.Range("D17").Value = '''0.798'
$ws.Range("E17").Value = '  +2.34%  '

$ws.Range("D18").Value = '42.187.91'
$ws.Range("E18").Value = '  +0.93%  '

$ws.Range("D19").Value = '''12.58'
$ws.Range("E19").Value = '  -2.41%  '

$ws.Range("D20").Value = '0.0₃0912'
$ws.Range("E20").Value = '  +0.19%  '

$ws.Range("D21").Value = '''6.00'
$ws.Range("E21").Value = '  +0.82%  '

$ws.Range("D22").Value = '''67.85'
$ws.Range("E22").Value = '  +0.85%  '

$ws.Range("D23").Value = '''238.21'
$ws.Range("E23").Value = '  -2.48%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '''2.59'
$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").Value = '''1.98'
$ws.Range("E25").Value = '  +2.33%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("D27").Value = '''23.81'
$ws.Range("E27").Value = '  -1.20%  '

$ws.Range("D28").Value = '''37.35'
$ws.Range("E28").Value = '  +6.11%  '

$ws.Range("E29").Value = '  -0.18%  '

$ws.Range("E30").Value = '  +1.85%  '

$ws.Range("D31").Value = '''160.65'
$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("D32").Value = '''5.28'
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").Value = '''3.18'
$ws.Range("E34").Value = '  +4.76%  '

$ws.Range("E35").Value = '  -0.33%  '

$ws.Range("D36").Value = '''17.33'
$ws.Range("E36").Value = '  +2.38%  '

$ws.Range("E37").Value = '  -1.05%  '

$ws.Range("E38").Value = '  -0.52%  '

$ws.Range("D39").Value = '''1.84'
$ws.Range("E39").Value = '  +1.80%  '

$ws.Range("E40").Value = '  -1.21%  '

$ws.Range("D41").Value = '''4.09'
$ws.Range("E41").Value = '  +3.48%  '

$ws.Range("E42").Value = '  +14.43%  '

$ws.Range("D43").Value = '1.989.22'
$ws.Range("E43").Value = '  -1.50%  '

$ws.Range("E44").Value = '  +1.36%  '

$ws.Range("D45").Value = '''18.95'
$ws.Range("E45").Value = '  -5.38%  '

$ws.Range("D46").Value = '''2.96'
$ws.Range("E46").Value = '  +1.58%  '

$ws.Range("D47").Value = '''9.98'
$ws.Range("E47").Value = '  -4.54%  '

$ws.Range("D48").Value = '''53.46'
$ws.Range("E48").Value = '  -0.13%  '

$ws.Range("E49").Value = '  +0.85%  '

$ws.Range("D50").Value = '''72.12'
$ws.Range("E50").Value = '  -1.45%  '

$ws.Range("D51").Value = '''91.72'
$ws.Range("E51").Value = '  -0.18%  '
